$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.437.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.609.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'211.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.0606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.834.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.611.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'63.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'235.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.26%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.420.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "  +4.53%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'6.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.492.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.561"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.747.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.759"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'61.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'89.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
$ws.Range("E51").Style = "Normal"
